# ranking_methodology.xlsx edit
# - Adds a new "Teaching quality" shared string (swallowed automatically by
#   the engine's shared-string table when we write the text) and re-labels
#   row 22 from "Research quality" to "Teaching quality" so it joins the
#   Family/Internal/Pediatrics/Surgery board-pass-rate rows underneath it.
# - Re-derives every Category-weight cell (column B) as
#   =SUMIFS($D$2:$D$27,$A$2:$A$27,A<row>) instead of a hard-coded literal.
# - Updates a batch of column-D metric weights.
# - Clears column D for the "Care quality" rows (13-19) -- no longer weighted.
# - Cosmetic: selection moves to D13, a 4th (D) column is introduced with a
#   width, and column B's width bumps very slightly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Column B: replace every literal weight with the SUMIFS formula ----
for ($r = 2; $r -le 27; $r++) {
    $ws.Range("B$r").Formula = '=SUMIFS($D$2:$D$27,$A$2:$A$27,$A' + $r + ')'
}

# ---- Column D: updated metric weights ----
$ws.Range("D2").Value = 0.2
$ws.Range("D7").Value = 0.06
$ws.Range("D8").Value = 0.02
$ws.Range("D9").Value = 0.02
$ws.Range("D10").Value = 0

# Care quality rows (13-19): weights cleared entirely (blank, style kept)
$ws.Range("D13").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("D15").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("D17").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("D19").ClearContents()

$ws.Range("D20").Value = 0.12
$ws.Range("D22").Value = 0.25
$ws.Range("D23").Value = 0.01
$ws.Range("D24").Value = 0.01
$ws.Range("D25").Value = 0.01
$ws.Range("D26").Value = 0.01
$ws.Range("D27").Value = 0.01

# ---- Row 22 moves from the "Research quality" category into "Teaching
# quality" (joining rows 23-27); C22 keeps its "Program setting" metric ----
$ws.Range("A22").Value = "Teaching quality"

# ---- New column D width + column B's tiny width bump ----
$ws.Columns.Item(2).ColumnWidth = 17.59
$ws.Columns.Item(4).ColumnWidth = 19.31

# ---- Selection / scroll position ----
$ws.Range("D13").Select()
